$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the task description in C4: replace the old UI/login design text
# with "Ghép nối các form" (form wiring task).
$ws.Range("C4").Value = "Ghép nối các form"

# Update the active selection to C4 to match the saved view state.
$ws.Activate()
$ws.Range("C4").Select()
